$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 9.24193
$ws.Range("N2").Value = 27.72579
$ws.Range("O2").Value = 0.1468938537243544
$ws.Range("P2").Value = 0.1569651396557324
$ws.Range("Q2").Value = 1.42870071677
$ws.Range("R2").Value = 12.85830645093
$ws.Range("S2").Value = 0.1468938537243544
$ws.Range("T2").Value = 0.1569651396557324

# Row 3
$ws.Range("O3").Value = 0.469548954544906
$ws.Range("P3").Value = 0.5017420086455576
$ws.Range("S3").Value = 0.469548954544906
$ws.Range("T3").Value = 0.5017420086455576

# Row 4
$ws.Range("M4").Value = 7.349831333333333
$ws.Range("N4").Value = 22.049494
$ws.Range("O4").Value = 0.1168203014713749
$ws.Range("P4").Value = 0.1248296948454213
$ws.Range("Q4").Value = 1.136203075988667
$ws.Range("R4").Value = 10.225827683898
$ws.Range("S4").Value = 0.1168203014713749
$ws.Range("T4").Value = 0.1248296948454213

# Row 5
$ws.Range("M5").Value = 12.1104985
$ws.Range("N5").Value = 24.220997
$ws.Range("O5").Value = 0.1924876941491673
$ws.Range("P5").Value = 0.1371233128688515
$ws.Range("Q5").Value = 1.8721498526165
$ws.Range("R5").Value = 11.232899115699
$ws.Range("S5").Value = 0.1924876941491673
$ws.Range("T5").Value = 0.1371233128688515

# Row 6
$ws.Range("M6").Value = 4.671440333333334
$ws.Range("N6").Value = 14.014321
$ws.Range("O6").Value = 0.07424919611019735
$ws.Range("P6").Value = 0.079339843984437
$ws.Range("Q6").Value = 0.7221532896896667
$ws.Range("R6").Value = 6.499379607207
$ws.Range("S6").Value = 0.07424919611019735
$ws.Range("T6").Value = 0.079339843984437
